$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.180.67"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "1.908.25"
$ws.Range("E3").Value = "  +5.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.95"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  +3.05%  "
$ws.Range("E8").Value = "  +4.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3021"
$ws.Range("E9").Value = "  +8.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06833"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").Value = "1.904.54"
$ws.Range("E11").Value = "  +5.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.35"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07321"
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6953"
$ws.Range("E14").Value = "  +7.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.31"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.937"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.165.02"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008307"
$ws.Range("E18").Value = "  +13.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9980"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.11"
$ws.Range("E20").Value = "  +6.53%  "
$ws.Range("D21").Value = "2.150.52"
$ws.Range("E21").Value = "  +5.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9973"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.835"
$ws.Range("E23").Value = "  +5.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.766"
$ws.Range("E24").Value = "  +7.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.361"
$ws.Range("E25").Value = "  +6.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.97"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.06"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.21"
$ws.Range("E28").Value = "  +4.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.012"
$ws.Range("E29").Value = "  +6.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.398"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.326"
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08890"
$ws.Range("E32").Value = "  +6.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.018"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05088"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.151"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7240"
$ws.Range("E36").Value = "  +7.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.312"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.821"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9642"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +6.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.081"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4328"
$ws.Range("E43").Value = "  +5.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.30"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.721"
$ws.Range("E46").Value = "  +7.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1283"
$ws.Range("E47").Value = "  +4.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05764"
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.45"
$ws.Range("E49").Value = "  +5.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.454"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3825"
$ws.Range("E51").Value = "  +5.47%  "
